# Fix typo 'Cyro-EM' to 'Cryo-EM' (Closes #9)
#
# The storage_medium list sheet re-sorted the "Cryo-EM" entry (after fixing
# its spelling) so it now appears earlier in the list - ahead of
# "DMSO (serum)" and "RNAlater", which both shift down by one row.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("storage_medium")

$ws.Range("A12").Value = "Cryo-EM"
$ws.Range("B12").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000333"

$ws.Range("A13").Value = "DMSO (serum)"
$ws.Range("B13").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000125"

$ws.Range("A14").Value = "RNAlater"
$ws.Range("B14").Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348"

# Update the .metadata sheet's pav:createdOn timestamp to reflect the new save.
$meta = $wb.Worksheets.Item(".metadata")
$meta.Range("C2").Value = "2024-03-14T10:55:17-04:00"
